$d = $word.ActiveDocument
$d.Content.Find.Execute("192.168.10.21", $true, $false, $false, $false, $false,
                         $true, 1, $false, "192.168.11.21", 2)
